$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 564, shifting rows 564:605 down to 565:606.
$ws.Rows(564).Insert()

# Populate the newly inserted row with the new data point for 2026/01/05 (Monday).
# Force column A to be written as text (matching the existing inlineStr/text cells
# in column A) instead of being auto-converted to a date serial number.
$ws.Range("A564").NumberFormat = "@"
$ws.Range("A564").Value = "2026/01/05"
$ws.Range("A564").ClearFormats()

$ws.Range("B564").Value = "月"
$ws.Range("C564").Value = 23
$ws.Range("D564").Value = 177
